# Apply the "pos-condicao" wording update for UC-37 (Cadastrar fornecedor)
# and move the automatic "_GoBack" bookmark to the end of that updated
# sentence (it previously sat, alone, inside an otherwise empty paragraph
# right after the first table in the document).

$d = $word.ActiveDocument

$oldSentence = "Fornecedor cadastrado com sucesso."
$newSentence = "Fornecedor cadastrado com sucesso e disponibilizado para o cadastro de manuten" + [char]0x00E7 + [char]0x00E3 + "o, m" + [char]0x00E1 + "quina e produtos."

# 1) Replace the short sentence with the expanded one. Word's Find/Replace
#    merges the (two) runs that used to make up the old sentence into a
#    single run carrying the new text, exactly like a user typing the
#    replacement in the Word UI would.
$replaced = $d.Content.Find.Execute($oldSentence, $true, $false, $false, $false, $false, `
                                     $true, 1, $false, $newSentence, 2)

if (-not $replaced) {
    throw "Could not find the target sentence to replace."
}

# 2) Re-locate the freshly inserted sentence so we know exactly where it
#    ends (this is where the relocated bookmark must collapse to).
$rng = $d.Content
$found = $rng.Find.Execute($newSentence, $true, $false, $false, $false, $false, `
                            $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not re-locate the replaced sentence."
}
$sentenceEnd = $rng.End

# 3) This engine's Bookmarks.Add refuses a truly collapsed range sitting on
#    a paragraph-boundary, so: insert a throw-away marker character right
#    after the sentence, wrap it with a non-collapsed bookmark range named
#    "_GoBack" (adding a bookmark with that name also removes/relocates any
#    pre-existing "_GoBack" bookmark elsewhere in the document), then erase
#    the marker character again. Deleting the bookmarked text collapses the
#    bookmark in place, leaving it exactly where the marker used to be -
#    i.e. right at the end of the new sentence, before the paragraph mark.
$insertPoint = $d.Range($sentenceEnd, $sentenceEnd)
$insertPoint.InsertAfter("X")

$markerRange = $d.Range($sentenceEnd, $sentenceEnd + 1)
$d.Bookmarks.Add("_GoBack", $markerRange)

$markerRange2 = $d.Range($sentenceEnd, $sentenceEnd + 1)
$markerRange2.Text = ""

Write-Output "Updated post-condition sentence and relocated the _GoBack bookmark."
